$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ46162189",
    "summ46273426",
    "summ46383773",
    "summ46494591",
    "summ46603445",
    "summ46722387",
    "summ46829347",
    "summ47052161",
    "summ47196111",
    "summ47328376",
    "summ47465727",
    "summ47601739",
    "summ47735837",
    "summ47879860",
    "summ48011176",
    "summ48124263",
    "summ48256455",
    "summ48396365",
    "summ48550062",
    "summ48714654",
    "summ48871935",
    "summ49067605",
    "summ49214975",
    "summ49368737",
    "summ49522743",
    "summ49672084",
    "summ49826940",
    "summ49996002",
    "summ50142404",
    "summ50282711",
    "summ50430313",
    "summ50582817",
    "summ50730385",
    "summ50882838",
    "summ51031040",
    "summ51190350",
    "summ51354629",
    "summ51511520",
    "summ51664366",
    "summ51832480",
    "summ51996758",
    "summ52142540",
    "summ52281180",
    "summ52418853",
    "summ52562052",
    "summ52706354",
    "summ52876457",
    "summ53025724",
    "summ53169527",
    "summ53331141"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

